# Insert the "Conferences and Seminars" schedule table.
#
# Strategy: Word COM has no direct "Tables.Add(rows-with-different-cell-counts)"
# primitive that reproduces an exact target OOXML shape (one table here has an
# irregular first row with a single spanning-width cell used purely as a
# coloured separator, followed by a normal 2-column data row). The reliable,
# fully-COM way to materialize an exact table shape is Range.InsertXML with a
# WordProcessingML "single xml part" package - the same mechanism Word itself
# uses under the hood for paste/InsertXML operations.
#
# 1) Append a new empty paragraph right after the heading paragraph - this
#    gives us an insertion host without touching the heading paragraph.
# 2) Call InsertXML on that new paragraph's Range; InsertXML *replaces* the
#    range's contents with the supplied WordOpenXML, so the placeholder
#    paragraph is swapped out for the table itself, leaving nothing behind.

$d = $word.ActiveDocument

$lastParaIndex = $d.Paragraphs.Count
$headingRange = $d.Paragraphs($lastParaIndex).Range
$headingRange.Collapse(0)
[void]$headingRange.InsertParagraphAfter()

$tableHost = $d.Paragraphs($lastParaIndex + 1).Range

$tableXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 wp14"><w:body><w:tbl><w:tblPr><w:tblCellMar><w:bottom w:type="auto" w:w="0"/><w:top w:type="auto" w:w="0"/><w:left w:type="auto" w:w="0"/><w:right w:type="auto" w:w="0"/></w:tblCellMar><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders><w:tblW w:type="auto" w:w="100"/></w:tblPr><w:tblGrid><w:gridCol w:w="100"/><w:gridCol w:w="100"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcBorders><w:top w:sz="0" w:color="white"/><w:bottom w:val="outset" w:sz="60" w:color="red"/><w:left w:sz="0" w:color="white"/><w:right w:sz="0" w:color="white"/></w:tcBorders><w:tcW w:w="12" w:type="pct"/></w:tcPr><w:p><w:r><w:t xml:space="preserve"></w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcBorders><w:top w:sz="0" w:color="white"/><w:bottom w:sz="0" w:color="white"/><w:left w:sz="0" w:color="white"/><w:right w:sz="0" w:color="white"/></w:tcBorders><w:tcW w:w="12" w:type="pct"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">2013-05-01</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcBorders><w:top w:sz="0" w:color="white"/><w:bottom w:sz="0" w:color="white"/><w:left w:sz="0" w:color="white"/><w:right w:sz="0" w:color="white"/></w:tcBorders><w:vAlign w:val="center"/><w:tcW w:w="88" w:type="pct"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">SemTechBiz San Francisco</w:t></w:r></w:p></w:tc></w:tr></w:tbl></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$tableHost.InsertXML($tableXml)

Write-Output "table inserted"
